$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "27.075.34"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "1.677.32"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "216.32"
$ws.Range("E5").Value = "  +1.53%  "
Set-TextValue "D6" "0.510"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "0.253"
$ws.Range("E8").Value = "  +3.13%  "
Set-TextValue "D9" "0.0619"
$ws.Range("E9").Value = "  +1.83%  "
Set-TextValue "D10" "20.20"
$ws.Range("E10").Value = "  +5.29%  "
Set-TextValue "D11" "0.0887"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").Value = "1.913.93"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").Value = "1.677.60"
$ws.Range("E13").Value = "  +3.61%  "
Set-TextValue "D14" "4.09"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.523"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "65.95"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.102.52"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D18" "238.40"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +1.84%  "
Set-TextValue "D20" "7.74"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("E23").Value = "  +2.93%  "
Set-TextValue "D24" "9.31"
Set-TextValue "D25" "145.74"
$ws.Range("E25").Value = "  -0.86%  "
Set-TextValue "D26" "7.14"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E27").Value = "  +0.55%  "
Set-TextValue "D28" "16.00"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("E29").Value = "  +0.01%  "
Set-TextValue "D30" "0.0499"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "1.477.09"
$ws.Range("E33").Value = "  -3.32%  "
Set-TextValue "D34" "3.11"
$ws.Range("E34").Value = "  +4.77%  "
$ws.Range("E35").Value = "  +6.15%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "0.576"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D38" "0.901"
$ws.Range("E38").Value = "  +8.61%  "
$ws.Range("E39").Value = "  +2.41%  "
Set-TextValue "D40" "6.07"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  -0.02%  "
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  +10.27%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "66.72"
$ws.Range("E43").Value = "  +8.48%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D44" "2.27"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").Value = "1.822.62"
$ws.Range("E45").Value = "  +3.82%  "
Set-TextValue "D46" "0.779"
$ws.Range("E46").Value = "  +2.42%  "
Set-TextValue "D47" "90.37"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("E49").Value = "  +4.98%  "
$ws.Range("E50").Value = "  +1.20%  "
Set-TextValue "D51" "7.69"
$ws.Range("E51").Value = "  +2.53%  "
